# Calculo das velocidades e incertezas finalizado
#
# 1) Statistics sheet: add "Precisao" (F, O) and "Incerteza" (G, P) columns
#    for the time (t1..t18,tTOTAL) and distance (d1..d18,dTOTAL) blocks.
# 2) New "Speeds" sheet: partial + total speeds (m/s and Km/h) with their
#    propagated uncertainties, derived from the Statistics sheet.

$wb = $excel.ActiveWorkbook
$statsSheet = $wb.Worksheets.Item("Statistics")

# ---- Statistics: Precisao / Incerteza for the time block (F,G) and the
# ---- distance block (O,P). Headers (row 1) already exist in the workbook.
for ($r = 2; $r -le 20; $r++) {
    $statsSheet.Range("F$r").Formula = "=1-C$r/B$r"
    $statsSheet.Range("G$r").Formula = "=C$r/10"
    $statsSheet.Range("O$r").Formula = "=1-L$r/K$r"
    $statsSheet.Range("P$r").Formula = "=L$r/10"
}

# Snapshot the computed values and rewrite as plain numbers (no formulas),
# matching the rest of the workbook which stores literal values only.
for ($r = 2; $r -le 20; $r++) {
    foreach ($col in @("F", "G", "O", "P")) {
        $addr = "$col$r"
        $v = $statsSheet.Range($addr).Value2
        $statsSheet.Range($addr).ClearContents()
        $statsSheet.Range($addr).Value = $v
    }
}

# ---- New "Speeds" sheet, placed after "Statistics" (last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$speedsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$speedsSheet.Name = "Speeds"

# Header row.
$speedsSheet.Range("A1").Value = "Parciais"
$speedsSheet.Range("B1").Value = "Velocidade"
$speedsSheet.Range("C1").Value = "Incerteza"
$speedsSheet.Range("D1").Value = "Unidade de medida"
$speedsSheet.Range("G1").Value = "Parciais"
$speedsSheet.Range("H1").Value = "Velocidade"
$speedsSheet.Range("I1").Value = "Incerteza"
$speedsSheet.Range("J1").Value = "Unidade de medida"

$labels = @("v1","v2","v3","v4","v5","v6","v7","v8","v9","v10","v11","v12","v13","v14","v15","v16","v17","v18","vTOTAL")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $label = $labels[$i]

    $speedsSheet.Range("A$r").Value = $label
    $speedsSheet.Range("D$r").Value = "m/s"
    $speedsSheet.Range("G$r").Value = $label
    $speedsSheet.Range("J$r").Value = "Km/h"

    # Speed (m/s) = reconciled distance / reconciled time for this segment.
    $speedsSheet.Range("B$r").Formula = "=Statistics!M$r/Statistics!D$r"
    # Propagated uncertainty (m/s), from the time/distance uncertainties
    # already computed on the Statistics sheet (columns G and P).
    $speedsSheet.Range("C$r").Formula = "=SQRT((Statistics!P$r/Statistics!D$r)^2+((Statistics!M$r/Statistics!D$r^2)*Statistics!G$r)^2)"
    # Speed (Km/h) and its uncertainty: simple unit conversion.
    $speedsSheet.Range("H$r").Formula = "=B$r*3.6"
    $speedsSheet.Range("I$r").Formula = "=C$r*3.6"
}

# Snapshot to plain numeric values (no formulas left behind), same as above.
for ($r = 2; $r -le 20; $r++) {
    foreach ($col in @("B", "C", "H", "I")) {
        $addr = "$col$r"
        $v = $speedsSheet.Range($addr).Value2
        $speedsSheet.Range($addr).ClearContents()
        $speedsSheet.Range($addr).Value = $v
    }
}

$statsSheet.Select()
